$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily-conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText.Replace("1000 Bs = 7.58 = 31003.79 pesos", "1000 Bs = 7.58 = 31012.27 pesos")
$newText = $newText.Replace("31003.79 pesos = 7.54 = 920.29 Bs", "31012.27 pesos = 7.56 = 961.42 Bs")
$cellA1.Value = $newText

# --- Sheet "tasas": update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 131.951
$ws2.Range("O10").Value = 4092.1
$ws2.Range("N12").Value = 4100
$ws2.Range("O12").Value = 127.105
